$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 542, shifting existing rows 542:642 down to 544:644
$ws.Rows("542:543").Insert()

# Populate new row 542 (Primera) with the new weekly data point
$ws.Range("A542").Value = 3
$ws.Range("B542").Value = "Femacal de La Calera"
$ws.Range("C542").Value = "Coquimbo"
$ws.Range("D542").Value = 44694
$ws.Range("E542").Value = 5
$ws.Range("F542").Value = 100112008
$ws.Range("G542").Value = "Coliflor"
$ws.Range("H542").Value = "Sin especificar"
$ws.Range("I542").Value = "Primera"
$ws.Range("J542").Value = 1830
$ws.Range("K542").Value = 1300
$ws.Range("L542").Value = 1500
$ws.Range("M542").Value = 1396
$ws.Range("N542").Value = "`$/unidad"
$ws.Range("O542").Value = "Provincia de Quillota"
$ws.Range("P542").Value = 1396
$ws.Range("Q542").Value = 1
$ws.Range("R542").Value = "Hortaliza"

# Populate new row 543 (Segunda) with the new weekly data point
$ws.Range("A543").Value = 3
$ws.Range("B543").Value = "Femacal de La Calera"
$ws.Range("C543").Value = "Coquimbo"
$ws.Range("D543").Value = 44694
$ws.Range("E543").Value = 5
$ws.Range("F543").Value = 100112008
$ws.Range("G543").Value = "Coliflor"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Segunda"
$ws.Range("J543").Value = 800
$ws.Range("K543").Value = 900
$ws.Range("L543").Value = 900
$ws.Range("M543").Value = 900
$ws.Range("N543").Value = "`$/unidad"
$ws.Range("O543").Value = "Provincia de Quillota"
$ws.Range("P543").Value = 900
$ws.Range("Q543").Value = 1
$ws.Range("R543").Value = "Hortaliza"
